$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so values like "30.459.40" or "1.005"
# are not auto-converted to numbers, matching the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.459.40"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "2.106.66"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "333.49"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").Value = "0.5242"
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("D8").Value = "0.4596"
$ws.Range("E8").Value = "  +5.59%  "
$ws.Range("D9").Value = "53.60"
$ws.Range("E9").Value = "  +13.15%  "
$ws.Range("D10").Value = "0.08951"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "1.176"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "24.37"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "2.097.69"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "6.771"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "7.843"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "96.49"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "0.00001128"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "0.06625"
$ws.Range("D20").Value = "19.23"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "6.282"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").Value = "30.535.49"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Value = "12.32"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").Value = "2.361"
$ws.Range("E25").Value = "  +3.24%  "
$ws.Range("D26").Value = "2.352.02"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "22.29"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").Value = "2.559"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").Value = "163.42"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").Value = "132.59"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "1.193"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "0.1071"
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").Value = "1.683"
$ws.Range("E33").Value = "  +8.83%  "
$ws.Range("D34").Value = "6.142"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").Value = "3.932"
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("D36").Value = "10.43"
$ws.Range("E36").Value = "  +8.52%  "
$ws.Range("D37").Value = "0.02567"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").Value = "0.06812"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "5.536"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "12.75"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").Value = "0.2288"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D42").Value = "0.6877"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "1.245"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "2.349"
$ws.Range("E44").Value = "  +5.87%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").Value = "0.6370"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "13.92"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").Value = "3.652"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").Value = "0.00000000350"
$ws.Range("E49").Value = "  +22.89%  "
$ws.Range("D50").Value = "1.244"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").Value = "1.220"
$ws.Range("E51").Value = "  +2.25%  "

# Restore default style on column D so no stray style index is left on the cells
$ws.Range("D2:D51").Style = "Normal"
